$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$arr = New-Object 'object[,]' 24,13
$arr[0,0] = 1.933140312091183
$arr[0,1] = 0.03175961529329641
$arr[0,2] = 0.5143914853724141
$arr[0,3] = 0.1593507817317494
$arr[0,4] = 0
$arr[0,5] = 3.019918314578348
$arr[0,6] = 2.357406099311703
$arr[0,7] = 0
$arr[0,8] = 0.05351804162523344
$arr[0,9] = 1.25908030626556
$arr[0,10] = 0.4885825854267409
$arr[0,11] = 0
$arr[0,12] = 3.752109414407087
$arr[1,0] = 1.896848963859497
$arr[1,1] = 0.02892879181674601
$arr[1,2] = 0.5122141626836054
$arr[1,3] = 0.1594562759811087
$arr[1,4] = 0
$arr[1,5] = 3.014944477446008
$arr[1,6] = 2.36102523265572
$arr[1,7] = 0
$arr[1,8] = 0.05367818756345111
$arr[1,9] = 1.223330531728209
$arr[1,10] = 0.4851973522642652
$arr[1,11] = 0
$arr[1,12] = 3.769650534517581
$arr[2,0] = 1.875576776482006
$arr[2,1] = 0.02717606729525102
$arr[2,2] = 0.5110903268864746
$arr[2,3] = 0.1595736432765982
$arr[2,4] = 0
$arr[2,5] = 3.01310038221277
$arr[2,6] = 2.364030361085838
$arr[2,7] = 0
$arr[2,8] = 0.05378294159167973
$arr[2,9] = 1.202087019899835
$arr[2,10] = 0.4833362672458605
$arr[2,11] = 0
$arr[2,12] = 3.781276122356843
$arr[3,0] = 1.867162717382854
$arr[3,1] = 0.02645812176991313
$arr[3,2] = 0.5106859997387971
$arr[3,3] = 0.1596347200088672
$arr[3,4] = 0
$arr[3,5] = 3.012653077550894
$arr[3,6] = 2.365451983845631
$arr[3,7] = 0
$arr[3,8] = 0.05382724908814795
$arr[3,9] = 1.193608072869466
$arr[3,10] = 0.4826326144821422
$arr[3,11] = 0
$arr[3,12] = 3.786228685145083
$arr[4,0] = 1.865780952568571
$arr[4,1] = 0.02633868341928292
$arr[4,2] = 0.5106221036992054
$arr[4,3] = 0.1596456626341496
$arr[4,4] = 0
$arr[4,5] = 3.012597173294367
$arr[4,6] = 2.365699945807762
$arr[4,7] = 0
$arr[4,8] = 0.05383470424065218
$arr[4,9] = 1.192210903679126
$arr[4,10] = 0.4825190826580013
$arr[4,11] = 0
$arr[4,12] = 3.787064041786856
$arr[5,0] = 1.87546227049873
$arr[5,1] = 0.02716639983257352
$arr[5,2] = 0.5110846566853695
$arr[5,3] = 0.1595744133040764
$arr[5,4] = 0
$arr[5,5] = 3.013093118169891
$arr[5,6] = 2.364048735750387
$arr[5,7] = 0
$arr[5,8] = 0.05378353257569213
$arr[5,9] = 1.201971949078995
$arr[5,10] = 0.4833265557427353
$arr[5,11] = 0
$arr[5,12] = 3.781342043701471
$arr[6,0] = 1.92041744228257
$arr[6,1] = 0.03078655742161374
$arr[6,2] = 0.513596567374961
$arr[6,3] = 0.1593762542173423
$arr[6,4] = 0
$arr[6,5] = 3.017952189377638
$arr[6,6] = 2.358491543089173
$arr[6,7] = 0
$arr[6,8] = 0.05357192913812892
$arr[6,9] = 1.246607137584903
$arr[6,10] = 0.4873702634958619
$arr[6,11] = 0
$arr[6,12] = 3.757980067228274
$arr[7,0] = 2.016587143999971
$arr[7,1] = 0.03777149204270813
$arr[7,2] = 0.5202105214063693
$arr[7,3] = 0.15940405484697
$arr[7,4] = 0
$arr[7,5] = 3.037086638540188
$arr[7,6] = 2.353802504281788
$arr[7,7] = 0
$arr[7,8] = 0.05320776231324409
$arr[7,9] = 1.339744933273494
$arr[7,10] = 0.497023258879878
$arr[7,11] = 0
$arr[7,12] = 3.718954648915442
$arr[8,0] = 2.092127973151207
$arr[8,1] = 0.04283631120298992
$arr[8,2] = 0.5260969777274482
$arr[8,3] = 0.1596772707495973
$arr[8,4] = 0
$arr[8,5] = 3.057014742726238
$arr[8,6] = 2.354139590013716
$arr[8,7] = 0
$arr[8,8] = 0.05297091798024489
$arr[8,9] = 1.411599958241027
$arr[8,10] = 0.5051644479355986
$arr[8,11] = 0
$arr[8,12] = 3.694422244015286
$arr[9,0] = 2.127555007406158
$arr[9,1] = 0.04512646243307472
$arr[9,2] = 0.5289976106311371
$arr[9,3] = 0.1598562477299588
$arr[9,4] = 0
$arr[9,5] = 3.06735866001182
$arr[9,6] = 2.355113832050279
$arr[9,7] = 0
$arr[9,8] = 0.05286978625530647
$arr[9,9] = 1.445035073100911
$arr[9,10] = 0.5090957072725928
$arr[9,11] = 0
$arr[9,12] = 3.684161198181954
$arr[10,0] = 2.141123034766679
$arr[10,1] = 0.04599172725592382
$arr[10,2] = 0.5301280030622593
$arr[10,3] = 0.1599318661422551
$arr[10,4] = 0
$arr[10,5] = 3.071459675115392
$arr[10,6] = 2.355600735139575
$arr[10,7] = 0
$arr[10,8] = 0.05283243671135196
$arr[10,9] = 1.457803627174428
$arr[10,10] = 0.5106170860882742
$arr[10,11] = 0
$arr[10,12] = 3.680404929764777
$arr[11,0] = 2.138194136036986
$arr[11,1] = 0.0458054639852179
$arr[11,2] = 0.5298831312067875
$arr[11,3] = 0.159915231838454
$arr[11,4] = 0
$arr[11,5] = 3.070568262348814
$arr[11,6] = 2.355490625712406
$arr[11,7] = 0
$arr[11,8] = 0.05284043855133724
$arr[11,9] = 1.455048916145046
$arr[11,10] = 0.5102879767082413
$arr[11,11] = 0
$arr[11,12] = 3.681208153928864
$arr[12,0] = 2.128668200974403
$arr[12,1] = 0.04519768753063147
$arr[12,2] = 0.529089967999667
$arr[12,3] = 0.159862311819591
$arr[12,4] = 0
$arr[12,5] = 3.067692364528284
$arr[12,6] = 2.355151525753456
$arr[12,7] = 0
$arr[12,8] = 0.05286669453092863
$arr[12,9] = 1.446083398057084
$arr[12,10] = 0.509220217081193
$arr[12,11] = 0
$arr[12,12] = 3.683849574362043
$arr[13,0] = 2.122853152140181
$arr[13,1] = 0.04482515171901014
$arr[13,2] = 0.5286082963368273
$arr[13,3] = 0.1598309176559916
$arr[13,4] = 0
$arr[13,5] = 3.065954760550994
$arr[13,6] = 2.354959179669351
$arr[13,7] = 0
$arr[13,8] = 0.05288290027613307
$arr[13,9] = 1.440605741983461
$arr[13,10] = 0.5085704397339725
$arr[13,11] = 0
$arr[13,12] = 3.685484371407455
$arr[14,0] = 2.089834100400026
$arr[14,1] = 0.04268636786284219
$arr[14,2] = 0.5259118937526068
$arr[14,3] = 0.1596666729134348
$arr[14,4] = 0
$arr[14,5] = 3.056364483082035
$arr[14,6] = 2.354092433915724
$arr[14,7] = 0
$arr[14,8] = 0.05297765987156833
$arr[14,9] = 1.409429931997067
$arr[14,10] = 0.50491210949221
$arr[14,11] = 0
$arr[14,12] = 3.695110924382192
$arr[15,0] = 2.069850060731255
$arr[15,1] = 0.04137076268398232
$arr[15,2] = 0.5243147731072355
$arr[15,3] = 0.1595799066963437
$arr[15,4] = 0
$arr[15,5] = 3.050808722440593
$arr[15,6] = 2.353770917792644
$arr[15,7] = 0
$arr[15,8] = 0.05303748216904935
$arr[15,9] = 1.390496032998726
$arr[15,10] = 0.5027261422627305
$arr[15,11] = 0
$arr[15,12] = 3.701246801921755
$arr[16,0] = 2.058455854159718
$arr[16,1] = 0.04061275513141993
$arr[16,2] = 0.5234171312239226
$arr[16,3] = 0.1595351499333688
$arr[16,4] = 0
$arr[16,5] = 3.047733527233447
$arr[16,6] = 2.353663260873475
$arr[16,7] = 0
$arr[16,8] = 0.05307251274892888
$arr[16,9] = 1.379676159709362
$arr[16,10] = 0.5014902760652404
$arr[16,11] = 0
$arr[16,12] = 3.704860595149654
$arr[17,0] = 2.054615170413058
$arr[17,1] = 0.04035588161688963
$arr[17,2] = 0.523116810242783
$arr[17,3] = 0.1595208810172295
$arr[17,4] = 0
$arr[17,5] = 3.046712981596755
$arr[17,6] = 2.353640083473863
$arr[17,7] = 0
$arr[17,8] = 0.05308448050747572
$arr[17,9] = 1.376024833118919
$arr[17,10] = 0.5010755179353623
$arr[17,11] = 0
$arr[17,12] = 3.706098690080665
$arr[18,0] = 2.071967039682818
$arr[18,1] = 0.04151094609605366
$arr[18,2] = 0.5244826184514721
$arr[18,3] = 0.1595886103696813
$arr[18,4] = 0
$arr[18,5] = 3.051387687735115
$arr[18,6] = 2.353797146986949
$arr[18,7] = 0
$arr[18,8] = 0.0530310495965427
$arr[18,9] = 1.392504293327278
$arr[18,10] = 0.5029566232365283
$arr[18,11] = 0
$arr[18,12] = 3.700584870824599
$arr[19,0] = 2.13146206043092
$arr[19,1] = 0.04537625933808442
$arr[19,2] = 0.5293220716550309
$arr[19,3] = 0.1598776429988575
$arr[19,4] = 0
$arr[19,5] = 3.068532090565753
$arr[19,6] = 2.355247926225587
$arr[19,7] = 0
$arr[19,8] = 0.05285895684492559
$arr[19,9] = 1.448713874333151
$arr[19,10] = 0.50953295707194
$arr[19,11] = 0
$arr[19,12] = 3.683070213263065
$arr[20,0] = 2.171234705098129
$arr[20,1] = 0.04789101992405165
$arr[20,2] = 0.5326713438605566
$arr[20,3] = 0.1601122523081919
$arr[20,4] = 0
$arr[20,5] = 3.080809503551649
$arr[20,6] = 2.356883737103885
$arr[20,7] = 0
$arr[20,8] = 0.05275200154545256
$arr[20,9] = 1.486076040169849
$arr[20,10] = 0.5140215230470062
$arr[20,11] = 0
$arr[20,12] = 3.672377476360651
$arr[21,0] = 2.149926030466418
$arr[21,1] = 0.04654988301487606
$arr[21,2] = 0.5308667374704044
$arr[21,3] = 0.1599828612007386
$arr[21,4] = 0
$arr[21,5] = 3.074158627957161
$arr[21,6] = 2.355947775432043
$arr[21,7] = 0
$arr[21,8] = 0.05280858197244953
$arr[21,9] = 1.466077921048878
$arr[21,10] = 0.5116084749166561
$arr[21,11] = 0
$arr[21,12] = 3.678015352897418
$arr[22,0] = 2.071009657917443
$arr[22,1] = 0.04144757431423329
$arr[22,2] = 0.5244066715134323
$arr[22,3] = 0.1595846594697825
$arr[22,4] = 0
$arr[22,5] = 3.051125567199222
$arr[22,6] = 2.353785048311124
$arr[22,7] = 0
$arr[22,8] = 0.05303395577540293
$arr[22,9] = 1.391596154864288
$arr[22,10] = 0.5028523577597781
$arr[22,11] = 0
$arr[22,12] = 3.700883861358477
$arr[23,0] = 1.989712878837025
$arr[23,1] = 0.03589383832527915
$arr[23,2] = 0.5182407515332983
$arr[23,3] = 0.1593520662263828
$arr[23,4] = 0
$arr[23,5] = 3.030880461574441
$arr[23,6] = 2.354406644242346
$arr[23,7] = 0
$arr[23,8] = 0.05330086816735591
$arr[23,9] = 1.313947504762012
$arr[23,10] = 0.4942275297567278
$arr[23,11] = 0
$arr[23,12] = 3.72878525569493
$ws.Range("B2:N25").Value = $arr
Write-Output "done"
